$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Date Awarded", "Contract #", "Qty", "Contract Total", "NSN", "Part Name", "Vendor Name", "P/N", "Preservation Method", "Due Date")
$data = @("3/5/2019", "65", "4654", "654", "464", "65464", "4", "4654", "654", "654")

for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $cell = $ws.Cells.Item(2, $i + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $data[$i]
    $cell.Style = "Normal"
}

$ws.Range("I3").Select()
